$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.45"
$ws.Range("E2").Value = "'-0.83%"
$ws.Range("D3").Value = "'27.19"
$ws.Range("E3").Value = "'-0.53%"
$ws.Range("D4").Value = "'4.655"
$ws.Range("E4").Value = "'-10.78%"
$ws.Range("D5").Value = "'0.05874"
$ws.Range("E5").Value = "'-0.98%"
$ws.Range("D6").Value = "'6.628"
$ws.Range("E6").Value = "'-1.20%"
$ws.Range("D7").Value = "'0.8585"
$ws.Range("E7").Value = "'-1.29%"
$ws.Range("D8").Value = "'0.9401"
$ws.Range("E8").Value = "'-5.83%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1403"
$ws.Range("E9").Value = "'-1.09%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.04239"
$ws.Range("E10").Value = "'19.33%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07096"
$ws.Range("E11").Value = "'-0.98%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("E12").Value = "'-0.39%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09151"
$ws.Range("E13").Value = "'-0.76%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001528"
$ws.Range("E14").Value = "'-0.83%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006072"
$ws.Range("E15").Value = "'0.18%"
$ws.Range("D16").Value = "'0.006231"
$ws.Range("E16").Value = "'6.38%"
$ws.Range("E17").Value = "'0.50%"
$ws.Range("D18").Value = "'3.204"
$ws.Range("E18").Value = "'-1.97%"
$ws.Range("E19").Value = "'-0.52%"
$ws.Range("D20").Value = "'0.3052"
$ws.Range("E20").Value = "'-2.89%"
$ws.Range("E21").Value = "'-0.41%"
$ws.Range("D22").Value = "'3.817"
$ws.Range("E22").Value = "'8.42%"
$ws.Range("D23").Value = "'0.04234"
$ws.Range("E23").Value = "'1.20%"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'0.08%"
$ws.Range("D25").Value = "'0.004284"
$ws.Range("E25").Value = "'-5.38%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.12%"
$ws.Range("E27").Value = "'0.13%"
$ws.Range("E40").Value = "'-0.37%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1101"
$ws.Range("E41").Value = "'-0.28%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003921"
$ws.Range("E42").Value = "'-40.36%"
$ws.Range("D43").Value = "'0.002431"
$ws.Range("E43").Value = "'3.09%"
$ws.Range("D44").Value = "'0.01140"
$ws.Range("E44").Value = "'6.51%"
$ws.Range("D45").Value = "'0.00005478"
$ws.Range("E45").Value = "'0.92%"
$ws.Range("E46").Value = "'0.18%"
$ws.Range("D47").Value = "'0.05002"
$ws.Range("E47").Value = "'-54.09%"
$ws.Range("D48").Value = "'0.2251"
$ws.Range("E48").Value = "'9,976.94%"
$ws.Range("E49").Value = "'0.18%"
$ws.Range("E50").Value = "'0.18%"
